$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.728.88"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "2.246.94"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'115.02"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'290.05"
$ws.Range("E6").Value = "  +9.54%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'46.48"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'9.14"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'0.896"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "2.582.04"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "2.253.90"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "42.819.43"
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  +10.83%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "'73.85"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'3.45"
$ws.Range("E22").Value = "  +20.66%  "
$ws.Range("D23").Value = "'2.36"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").Value = "'232.84"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'9.30"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  +5.58%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'40.24"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("B30").Value = "WEMIXToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D30").Value = "'3.28"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'175.49"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.31"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0915"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.56"
$ws.Range("E34").Value = "  +18.38%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.65"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.129"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.71"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0373"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.106"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.66"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'13.59"
$ws.Range("E41").Value = "  -6.07%  "
$ws.Range("D42").Value = "'72.61"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.238"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.35"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "'5.60"
$ws.Range("E46").Value = "  -6.26%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'107.67"
$ws.Range("E47").Value = "  +7.23%  "
$ws.Range("D48").Value = "'1.32"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("D49").Value = "'8.60"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.652"
$ws.Range("E50").Value = "  +6.30%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.475"
$ws.Range("E51").Value = "  +8.96%  "
